$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "longname" column (F) for rows 5 and 6
$ws.Range("F5").Value = "Transport Canada Dash 7 - CGCFR"
$ws.Range("F6").Value = "Transport Canada Dash 8 - CGCFJ"

# Apply font formatting: Arial, size 12, color FF212529
$fontRange = $ws.Range("F5:F6")
$fontRange.Font.Name = "Arial"
$fontRange.Font.Size = 12
$fontRange.Font.Color = 2696481

# Set row heights for rows 5 and 6
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75

# Update the selected cell to F6
$ws.Range("F6").Select()
